$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before D; existing D:K data shifts to F:M
$ws.Columns("D:E").Insert()

# Copy number formatting (incl. date format row) into the two new columns
# for each contiguous block of data rows (skipping section-label rows 5,6,37,79
# and the blank separator rows 36,78 which must NOT receive new cells)
$ws.Range("F7:G35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:G77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:G102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Write the final values for every data cell D:M across all data rows
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("F7").Value = 43281
$ws.Range("G7").Value = 43190
$ws.Range("H7").Value = 43100
$ws.Range("I7").Value = 43008
$ws.Range("J7").Value = 42916
$ws.Range("K7").Value = 42825
$ws.Range("L7").Value = 42735
$ws.Range("M7").Value = 42643
$ws.Range("D8").Value = 309900
$ws.Range("E8").Value = 279400
$ws.Range("F8").Value = 258400
$ws.Range("G8").Value = 235200
$ws.Range("H8").Value = 329600
$ws.Range("I8").Value = 266200
$ws.Range("J8").Value = 278100
$ws.Range("K8").Value = 363000
$ws.Range("L8").Value = 410200
$ws.Range("M8").Value = 385200
$ws.Range("D9").Value = 193400
$ws.Range("E9").Value = 172700
$ws.Range("F9").Value = 159700
$ws.Range("G9").Value = 141200
$ws.Range("H9").Value = 158400
$ws.Range("I9").Value = 169900
$ws.Range("J9").Value = 167200
$ws.Range("K9").Value = 165900
$ws.Range("L9").Value = 182900
$ws.Range("M9").Value = 216300
$ws.Range("D10").Value = 116500
$ws.Range("E10").Value = 106700
$ws.Range("F10").Value = 98700
$ws.Range("G10").Value = 94000
$ws.Range("H10").Value = 171200
$ws.Range("I10").Value = 96300
$ws.Range("J10").Value = 110900
$ws.Range("K10").Value = 197100
$ws.Range("L10").Value = 227300
$ws.Range("M10").Value = 168900
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("F12").Value = "NA"
$ws.Range("G12").Value = "NA"
$ws.Range("H12").Value = "NA"
$ws.Range("I12").Value = "NA"
$ws.Range("J12").Value = "NA"
$ws.Range("K12").Value = "NA"
$ws.Range("L12").Value = "NA"
$ws.Range("M12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 0
$ws.Range("D14").Value = 2400
$ws.Range("E14").Value = -100
$ws.Range("F14").Value = 792800
$ws.Range("G14").Value = 8800
$ws.Range("H14").Value = 121600
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 1435400
$ws.Range("M14").Value = 0
$ws.Range("D15").Value = 114200
$ws.Range("E15").Value = 113900
$ws.Range("F15").Value = 129700
$ws.Range("G15").Value = 128800
$ws.Range("H15").Value = 138100
$ws.Range("I15").Value = 137600
$ws.Range("J15").Value = 136600
$ws.Range("K15").Value = 135700
$ws.Range("L15").Value = 155200
$ws.Range("M15").Value = 155200
$ws.Range("D17").Value = 324800
$ws.Range("E17").Value = 301100
$ws.Range("F17").Value = 1104000
$ws.Range("G17").Value = 300800
$ws.Range("H17").Value = 439900
$ws.Range("I17").Value = 322800
$ws.Range("J17").Value = 322400
$ws.Range("K17").Value = 317500
$ws.Range("L17").Value = 1788300
$ws.Range("M17").Value = 387400
$ws.Range("D18").Value = -14900
$ws.Range("E18").Value = -21700
$ws.Range("F18").Value = -845600
$ws.Range("G18").Value = -65600
$ws.Range("H18").Value = -110300
$ws.Range("I18").Value = -56600
$ws.Range("J18").Value = -44300
$ws.Range("K18").Value = 45500
$ws.Range("L18").Value = -1378100
$ws.Range("M18").Value = -2200
$ws.Range("D20").Value = 1500
$ws.Range("E20").Value = 2600
$ws.Range("F20").Value = 2800
$ws.Range("G20").Value = 1300
$ws.Range("H20").Value = 1800
$ws.Range("I20").Value = 1400
$ws.Range("J20").Value = 3100
$ws.Range("K20").Value = 1600
$ws.Range("L20").Value = 1400
$ws.Range("M20").Value = 500
$ws.Range("D21").Value = 100800
$ws.Range("E21").Value = 94700
$ws.Range("F21").Value = -713100
$ws.Range("G21").Value = 64400
$ws.Range("H21").Value = 29600
$ws.Range("I21").Value = 82400
$ws.Range("J21").Value = 95400
$ws.Range("K21").Value = 182800
$ws.Range("L21").Value = -1221500
$ws.Range("M21").Value = 153600
$ws.Range("D22").Value = 73700
$ws.Range("E22").Value = 73700
$ws.Range("F22").Value = 74100
$ws.Range("G22").Value = 76000
$ws.Range("H22").Value = 72400
$ws.Range("I22").Value = 72900
$ws.Range("J22").Value = 73200
$ws.Range("K22").Value = 73400
$ws.Range("L22").Value = 55900
$ws.Range("M22").Value = 52600
$ws.Range("D23").Value = -87100
$ws.Range("E23").Value = -92800
$ws.Range("F23").Value = -916900
$ws.Range("G23").Value = -140300
$ws.Range("H23").Value = -180900
$ws.Range("I23").Value = -128100
$ws.Range("J23").Value = -114400
$ws.Range("K23").Value = -26400
$ws.Range("L23").Value = -1432600
$ws.Range("M23").Value = -54200
$ws.Range("D24").Value = -56300
$ws.Range("E24").Value = -14500
$ws.Range("F24").Value = -38800
$ws.Range("G24").Value = 3000
$ws.Range("H24").Value = -59000
$ws.Range("I24").Value = -28600
$ws.Range("J24").Value = -18200
$ws.Range("K24").Value = 257400
$ws.Range("L24").Value = -149500
$ws.Range("M24").Value = -10000
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 0
$ws.Range("D26").Value = -30800
$ws.Range("E26").Value = -78400
$ws.Range("F26").Value = -878000
$ws.Range("G26").Value = -143300
$ws.Range("H26").Value = -122000
$ws.Range("I26").Value = -99500
$ws.Range("J26").Value = -96200
$ws.Range("K26").Value = -283800
$ws.Range("L26").Value = -1283200
$ws.Range("M26").Value = -44200
$ws.Range("D27").Value = -33100
$ws.Range("E27").Value = -81600
$ws.Range("F27").Value = -628100
$ws.Range("G27").Value = -142300
$ws.Range("H27").Value = -133700
$ws.Range("I27").Value = -96800
$ws.Range("J27").Value = -91900
$ws.Range("K27").Value = -301700
$ws.Range("L27").Value = -1289400
$ws.Range("M27").Value = -55100
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = "NA"
$ws.Range("H29").Value = 109000
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = -1500
$ws.Range("K29").Value = "NA"
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = 0
$ws.Range("D32").Value = -1500
$ws.Range("E32").Value = -2600
$ws.Range("F32").Value = -2800
$ws.Range("G32").Value = -1300
$ws.Range("H32").Value = -1800
$ws.Range("I32").Value = -1400
$ws.Range("J32").Value = -3100
$ws.Range("K32").Value = -1600
$ws.Range("L32").Value = -1400
$ws.Range("M32").Value = -500
$ws.Range("D33").Value = -33100
$ws.Range("E33").Value = -81600
$ws.Range("F33").Value = -628100
$ws.Range("G33").Value = -142300
$ws.Range("H33").Value = -24700
$ws.Range("I33").Value = -96800
$ws.Range("J33").Value = -93400
$ws.Range("K33").Value = -301700
$ws.Range("L33").Value = -1289400
$ws.Range("M33").Value = -55100
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = 0
$ws.Range("D35").Value = -33100
$ws.Range("E35").Value = -81600
$ws.Range("F35").Value = -628100
$ws.Range("G35").Value = -142300
$ws.Range("H35").Value = -24700
$ws.Range("I35").Value = -96800
$ws.Range("J35").Value = -93400
$ws.Range("K35").Value = -301700
$ws.Range("L35").Value = -1289400
$ws.Range("M35").Value = -55100
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("F38").Value = 43281
$ws.Range("G38").Value = 43190
$ws.Range("H38").Value = 43100
$ws.Range("I38").Value = 43008
$ws.Range("J38").Value = 42916
$ws.Range("K38").Value = 42825
$ws.Range("L38").Value = 42735
$ws.Range("M38").Value = 42643
$ws.Range("D41").Value = 375200
$ws.Range("E41").Value = 326300
$ws.Range("F41").Value = 411500
$ws.Range("G41").Value = 461700
$ws.Range("H41").Value = 662800
$ws.Range("I41").Value = 608800
$ws.Range("J41").Value = 603000
$ws.Range("K41").Value = 519800
$ws.Range("L41").Value = 725700
$ws.Range("M41").Value = 426100
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("F42").Value = 0
$ws.Range("G42").Value = 0
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = 0
$ws.Range("D43").Value = 221200
$ws.Range("E43").Value = 264700
$ws.Range("F43").Value = 225000
$ws.Range("G43").Value = 203300
$ws.Range("H43").Value = 310000
$ws.Range("I43").Value = 257900
$ws.Range("J43").Value = 260800
$ws.Range("K43").Value = 381800
$ws.Range("L43").Value = 374600
$ws.Range("M43").Value = 355000
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = 0
$ws.Range("D45").Value = 62600
$ws.Range("E45").Value = 69800
$ws.Range("F45").Value = 60700
$ws.Range("G45").Value = 55400
$ws.Range("H45").Value = 66100
$ws.Range("I45").Value = 74600
$ws.Range("J45").Value = 74300
$ws.Range("K45").Value = 63000
$ws.Range("L45").Value = 92300
$ws.Range("M45").Value = 102800
$ws.Range("D46").Value = 659100
$ws.Range("E46").Value = 660800
$ws.Range("F46").Value = 697300
$ws.Range("G46").Value = 720500
$ws.Range("H46").Value = 1039000
$ws.Range("I46").Value = 941300
$ws.Range("J46").Value = 938100
$ws.Range("K46").Value = 964500
$ws.Range("L46").Value = 1192600
$ws.Range("M46").Value = 883800
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("F47").Value = 0
$ws.Range("G47").Value = 0
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = 0
$ws.Range("D48").Value = 8480700
$ws.Range("E48").Value = 8543400
$ws.Range("F48").Value = 8521400
$ws.Range("G48").Value = 9398900
$ws.Range("H48").Value = 9489200
$ws.Range("I48").Value = 9712300
$ws.Range("J48").Value = 9838300
$ws.Range("K48").Value = 9944400
$ws.Range("L48").Value = 10061900
$ws.Range("M48").Value = 11591800
$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 0
$ws.Range("F49").Value = 0
$ws.Range("G49").Value = 0
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = 0
$ws.Range("D52").Value = 125100
$ws.Range("E52").Value = 196900
$ws.Range("F52").Value = 175000
$ws.Range("G52").Value = 148800
$ws.Range("H52").Value = 266400
$ws.Range("I52").Value = 244700
$ws.Range("J52").Value = 248700
$ws.Range("K52").Value = 97100
$ws.Range("L52").Value = 185600
$ws.Range("M52").Value = 108600
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = 0
$ws.Range("D54").Value = 9264900
$ws.Range("E54").Value = 9401100
$ws.Range("F54").Value = 9393700
$ws.Range("G54").Value = 10268100
$ws.Range("H54").Value = 10794700
$ws.Range("I54").Value = 10898200
$ws.Range("J54").Value = 11025100
$ws.Range("K54").Value = 11006000
$ws.Range("L54").Value = 11440100
$ws.Range("M54").Value = 12584100
$ws.Range("D57").Value = 125600
$ws.Range("E57").Value = 103300
$ws.Range("F57").Value = 93600
$ws.Range("G57").Value = 94300
$ws.Range("H57").Value = 84000
$ws.Range("I57").Value = 84000
$ws.Range("J57").Value = 86600
$ws.Range("K57").Value = 83800
$ws.Range("L57").Value = 108200
$ws.Range("M57").Value = 114400
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("F58").Value = 0
$ws.Range("G58").Value = 0
$ws.Range("H58").Value = 249800
$ws.Range("I58").Value = 249700
$ws.Range("J58").Value = 249500
$ws.Range("K58").Value = 249300
$ws.Range("L58").Value = 299900
$ws.Range("M58").Value = 299800
$ws.Range("D59").Value = 239900
$ws.Range("E59").Value = 212500
$ws.Range("F59").Value = 243600
$ws.Range("G59").Value = 200200
$ws.Range("H59").Value = 259100
$ws.Range("I59").Value = 261600
$ws.Range("J59").Value = 312300
$ws.Range("K59").Value = 216800
$ws.Range("L59").Value = 225200
$ws.Range("M59").Value = 270000
$ws.Range("D60").Value = 365500
$ws.Range("E60").Value = 315800
$ws.Range("F60").Value = 337200
$ws.Range("G60").Value = 294500
$ws.Range("H60").Value = 593000
$ws.Range("I60").Value = 595300
$ws.Range("J60").Value = 648500
$ws.Range("K60").Value = 549900
$ws.Range("L60").Value = 633300
$ws.Range("M60").Value = 684100
$ws.Range("D61").Value = 3877400
$ws.Range("E61").Value = 3903000
$ws.Range("F61").Value = 3842600
$ws.Range("G61").Value = 3841400
$ws.Range("H61").Value = 3795900
$ws.Range("I61").Value = 3795300
$ws.Range("J61").Value = 3793900
$ws.Range("K61").Value = 3792500
$ws.Range("L61").Value = 4040200
$ws.Range("M61").Value = 3830200
$ws.Range("D62").Value = 367500
$ws.Range("E62").Value = 490600
$ws.Range("F62").Value = 440800
$ws.Range("G62").Value = 473500
$ws.Range("H62").Value = 455100
$ws.Range("I62").Value = 542800
$ws.Range("J62").Value = 510300
$ws.Range("K62").Value = 481700
$ws.Range("L62").Value = 299200
$ws.Range("M62").Value = 311800
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = 0
$ws.Range("D66").Value = 5011800
$ws.Range("E66").Value = 5113100
$ws.Range("F66").Value = 5031500
$ws.Range("G66").Value = 5280200
$ws.Range("H66").Value = 5518500
$ws.Range("I66").Value = 5611300
$ws.Range("J66").Value = 5648700
$ws.Range("K66").Value = 5545400
$ws.Range("L66").Value = 5681400
$ws.Range("M66").Value = 5539200
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = 0
$ws.Range("D72").Value = 3608400
$ws.Range("E72").Value = 3641400
$ws.Range("F72").Value = 3721400
$ws.Range("G72").Value = 4351100
$ws.Range("H72").Value = 4637700
$ws.Range("I72").Value = 4662500
$ws.Range("J72").Value = 4759300
$ws.Range("K72").Value = 4852600
$ws.Range("L72").Value = 5154200
$ws.Range("M72").Value = 6457100
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = 0
$ws.Range("D76").Value = 4253200
$ws.Range("E76").Value = 4288000
$ws.Range("F76").Value = 4362200
$ws.Range("G76").Value = 4988000
$ws.Range("H76").Value = 5276200
$ws.Range("I76").Value = 5287000
$ws.Range("J76").Value = 5376400
$ws.Range("K76").Value = 5460500
$ws.Range("L76").Value = 5758700
$ws.Range("M76").Value = 7044900
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("F80").Value = 43281
$ws.Range("G80").Value = 43190
$ws.Range("H80").Value = 43100
$ws.Range("I80").Value = 43008
$ws.Range("J80").Value = 42916
$ws.Range("K80").Value = 42825
$ws.Range("L80").Value = 42735
$ws.Range("M80").Value = 42643
$ws.Range("D81").Value = -33100
$ws.Range("E81").Value = -81600
$ws.Range("F81").Value = -628100
$ws.Range("G81").Value = -142300
$ws.Range("H81").Value = -24700
$ws.Range("I81").Value = -96800
$ws.Range("J81").Value = -93400
$ws.Range("K81").Value = -301700
$ws.Range("L81").Value = -1289400
$ws.Range("M81").Value = -55100
$ws.Range("D83").Value = 114200
$ws.Range("E83").Value = 113900
$ws.Range("F83").Value = 129700
$ws.Range("G83").Value = 128800
$ws.Range("H83").Value = 138100
$ws.Range("I83").Value = 137600
$ws.Range("J83").Value = 136600
$ws.Range("K83").Value = 135700
$ws.Range("L83").Value = 155200
$ws.Range("M83").Value = 155200
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("F84").Value = 0
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("F85").Value = 0
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("F86").Value = 0
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("F87").Value = 0
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = 0
$ws.Range("D89").Value = 128500
$ws.Range("E89").Value = -10600
$ws.Range("F89").Value = -1000
$ws.Range("G89").Value = 54900
$ws.Range("H89").Value = 154900
$ws.Range("I89").Value = 44800
$ws.Range("J89").Value = 112400
$ws.Range("K89").Value = 141900
$ws.Range("L89").Value = 165600
$ws.Range("M89").Value = 101100
$ws.Range("D91").Value = -45500
$ws.Range("E91").Value = -73500
$ws.Range("F91").Value = -42100
$ws.Range("G91").Value = -33800
$ws.Range("H91").Value = -36800
$ws.Range("I91").Value = -25400
$ws.Range("J91").Value = -30200
$ws.Range("K91").Value = -18700
$ws.Range("L91").Value = -727800
$ws.Range("M91").Value = -351000
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 0
$ws.Range("D94").Value = -44200
$ws.Range("E94").Value = -73100
$ws.Range("F94").Value = -38400
$ws.Range("G94").Value = -33700
$ws.Range("H94").Value = -32900
$ws.Range("I94").Value = -18100
$ws.Range("J94").Value = -29200
$ws.Range("K94").Value = -38100
$ws.Range("L94").Value = -60000
$ws.Range("M94").Value = -472200
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("F96").Value = 0
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -5000
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 0
$ws.Range("D100").Value = -35400
$ws.Range("E100").Value = -800
$ws.Range("F100").Value = -10700
$ws.Range("G100").Value = -222400
$ws.Range("H100").Value = -30600
$ws.Range("I100").Value = -20900
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = -309700
$ws.Range("L100").Value = 196300
$ws.Range("M100").Value = -28200
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("F101").Value = 0
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = 0
$ws.Range("D102").Value = 48900
$ws.Range("E102").Value = -84500
$ws.Range("F102").Value = -50200
$ws.Range("G102").Value = -201200
$ws.Range("H102").Value = 54100
$ws.Range("I102").Value = 5800
$ws.Range("J102").Value = 83200
$ws.Range("K102").Value = -206000
$ws.Range("L102").Value = 299700
$ws.Range("M102").Value = -396900
